# Add basic report functions:
# - Record the playing field (Spielfeld) and game class (Spielklasse) for game No. 1
# - Clear the (previously hard-coded) team name placeholders so the sheet
#   relies on the formulas already in place (L2/P2 reference B22/F22)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D7").Value = 1
$ws.Range("E7").Value = "Mixed"

$ws.Range("B22").Value = ""
$ws.Range("F22").Value = ""

[void]$ws.Range("X28").Select()
